$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 51; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = $cell.Value2 * 100
}

$ws.Range("B2:B51").NumberFormat = "General"
